$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 8.8626

$ws.Range("B7").Value = 4.872299999999997
$ws.Range("C7").Value = -13.96329999999999

$ws.Range("C15").Value = -14.61079999999999

$ws.Range("B16").Value = 7.502699999999996
$ws.Range("D16").Value = -8.636500000000003

$ws.Range("D19").Value = -8.160999999999996

$ws.Range("C21").Value = -12.3082

$ws.Range("C22").Value = -12.79540000000001

$ws.Range("C23").Value = -12.3962

$ws.Range("B28").Value = 5.762299999999999

$ws.Range("B29").Value = 5.182900000000002

$ws.Range("B32").Value = 7.193599999999995

$ws.Range("C34").Value = -11.37960000000002

$ws.Range("D36").Value = -8.543799999999994

$ws.Range("B40").Value = 9.017099999999994

$ws.Range("C43").Value = -13.08549999999999

$ws.Range("C45").Value = -13.66909999999999

$ws.Range("D46").Value = -8.3779

$ws.Range("C50").Value = -13.86419999999999
$ws.Range("D50").Value = -8.148400000000002

$ws.Range("C51").Value = -11.7305

$ws.Range("B52").Value = 5.3568

$ws.Range("B57").Value = 5.013899999999998

$ws.Range("B66").Value = 5.965300000000003
$ws.Range("C66").Value = -11.46630000000001

$ws.Range("C67").Value = -11.247

$ws.Range("C79").Value = -11.4346

$ws.Range("C84").Value = -12.6186

$ws.Range("C92").Value = -11.4385

$ws.Range("D95").Value = -8.0891

$ws.Range("C97").Value = -11.9066
$ws.Range("D97").Value = -8.680599999999998

$ws.Range("B100").Value = 5.915599999999999
